# Generate Report for Handoff
# Status moves from "In Translation" to "Ready for handoff" and the
# "Latest Handoff Datetime" / summary timestamps advance to the new
# handoff-generation time. Re-touching the Status column also makes the
# report widen that column slightly to fit the new, longer label.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-10-14 08:24:30"
$overview.Columns.Item(5).ColumnWidth = 16.28
$overview.Columns.Item(6).ColumnWidth = 16.28

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-10-14 08:24:20"
$zhcn.Columns.Item(3).ColumnWidth = 16.28

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-10-14 08:24:30"
$dede.Columns.Item(3).ColumnWidth = 16.28
